$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1): shift header labels C1->prediction, D1->rejection-f, E1->max
$ws.Range("C1").Value = "prediction"
$ws.Range("D1").Value = "rejection-f"
$ws.Range("E1").Value = "max"

# Data rows 2-9: column C becomes the text label "f__CAG-313",
# column D stays "f__CAG-313", column E becomes numeric 1
for ($r = 2; $r -le 9; $r++) {
    $ws.Cells.Item($r, 3).Value = "f__CAG-313"
    $ws.Cells.Item($r, 4).Value = "f__CAG-313"
    $ws.Cells.Item($r, 5).Value = 1
}
